# Applies the diff to the betexplorer "slovenia 2-snl 2023-2024" sheet:
#   1) Several existing rows (3,4,6 / 13,15,16 / 29,30 / 33,34,35,36 / 46,47)
#      have their match data (columns F..V) rotated among themselves - the
#      row identity / Indice (col A) and date (col E) stay put, only the
#      match details move.
#   2) Four brand-new match rows are appended at the end (58..61).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: read columns F..V (as an array, index 0 => col F) from a given row
# (NB: this COM shim only binds function args positionally, not by -name)
# ---------------------------------------------------------------------------
function Get-RowFV($row) {
    $vals = @()
    for ($c = 6; $c -le 22; $c++) {
        $vals += $ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowFV($row, $vals) {
    for ($c = 6; $c -le 22; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$c - 6]
    }
}

# ---------------------------------------------------------------------------
# 1) Snapshot the F..V data for every row that is involved in a reshuffle,
#    BEFORE writing anything back (several rows are both a source and a
#    destination, so all reads must happen first).
# ---------------------------------------------------------------------------
$snap = @{}
foreach ($r in 3,4,6,13,15,16,29,30,33,34,35,36,46,47) {
    $snap[$r] = Get-RowFV $r
}

# new_row -> old_row (source of the F..V data that ends up there)
$rowSource = @{
    3  = 6
    4  = 3
    6  = 4
    13 = 16
    15 = 13
    16 = 15
    29 = 30
    30 = 29
    33 = 34
    34 = 35
    35 = 36
    36 = 33
    46 = 47
    47 = 46
}

foreach ($dest in $rowSource.Keys) {
    $src = $rowSource[$dest]
    Set-RowFV $dest $snap[$src]
}

# ---------------------------------------------------------------------------
# 2) Append the four new rows (58..61), copying the row-57 formatting first
#    so the new cells inherit the same styles (bold/bordered Indice column,
#    date-formatted data_partida column) without creating new style entries.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row=58; A=57; E=45189.66666666666; F="NK Bistrica"; G=2;  H="Fuzinar";     I=4;
       J=1.45; K="19/09/2023 03:12"; L=1.47; M="20/09/2023 15:57";
       N=4.06; O="19/09/2023 03:12"; P=4.44; Q="20/09/2023 15:57";
       R=5.01; S="19/09/2023 03:12"; T=5.58; U="20/09/2023 15:57";
       V="https://www.betexplorer.com/football/slovenia/2-snl/bistrica-fuzinar/2ZRp5MNe/" },

    @{ Row=59; A=58; E=45189.66666666666; F="Dravinja"; G=5; H="Tabor Sezana"; I=2;
       J=1.77; K="19/09/2023 03:12"; L=2.2;  M="20/09/2023 15:50";
       N=3.5;  O="19/09/2023 03:12"; P=3.33; Q="20/09/2023 15:52";
       R=3.6;  S="19/09/2023 03:12"; T=3.02; U="20/09/2023 15:50";
       V="https://www.betexplorer.com/football/slovenia/2-snl/dravinja-tabor-sezana/f9Ox7rgq/" },

    @{ Row=60; A=59; E=45189.66666666666; F="Ilirija"; G=0; H="ND Gorica"; I=4;
       J=2.65; K="19/09/2023 03:12"; L=4.18; M="20/09/2023 15:57";
       N=3.15; O="19/09/2023 03:12"; P=4.37; Q="20/09/2023 15:57";
       R=2.31; S="19/09/2023 03:12"; T=1.62; U="20/09/2023 15:57";
       V="https://www.betexplorer.com/football/slovenia/2-snl/ilirija-nd-gorica/tQQl4tw2/" },

    @{ Row=61; A=60; E=45189.75; F="NK Krka"; G=2; H="Primorje"; I=3;
       J=2.33; K="19/09/2023 05:12"; L=3.87; M="20/09/2023 17:32";
       N=3.14; O="19/09/2023 05:12"; P=3.37; Q="20/09/2023 17:35";
       R=2.64; S="19/09/2023 05:12"; T=1.81; U="20/09/2023 17:32";
       V="https://www.betexplorer.com/football/slovenia/2-snl/nk-krka-primorje/8CSt628k/" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    # Clone formatting from the last existing data row (57) so styles
    # (bold+border Indice cell, date-formatted data_partida cell) line up
    # with the existing sheet's style table instead of minting new ones.
    $ws.Range("A57:V57").Copy()
    $ws.Range("A" + $r + ":V" + $r).PasteSpecial(-4122)   # xlPasteFormats

    $ws.Cells.Item($r, 1).Value  = $nr.A           # Indice
    $ws.Cells.Item($r, 2).Value  = "slovenia"      # pais
    $ws.Cells.Item($r, 3).Value  = "2-snl"         # torneio
    $ws.Cells.Item($r, 4).Value  = "2023-2024"     # temporada
    $ws.Cells.Item($r, 5).Value  = $nr.E           # data_partida
    $ws.Cells.Item($r, 6).Value  = $nr.F           # home
    $ws.Cells.Item($r, 7).Value  = $nr.G           # home_ft_gols
    $ws.Cells.Item($r, 8).Value  = $nr.H           # away
    $ws.Cells.Item($r, 9).Value  = $nr.I           # away_ft_gols
    $ws.Cells.Item($r, 10).Value = $nr.J           # home_opening_odds
    $ws.Cells.Item($r, 11).Value = $nr.K           # home_opening_data_hora
    $ws.Cells.Item($r, 12).Value = $nr.L           # home_closing_odds
    $ws.Cells.Item($r, 13).Value = $nr.M           # home_closing_data_hora
    $ws.Cells.Item($r, 14).Value = $nr.N           # draw_opening_odds
    $ws.Cells.Item($r, 15).Value = $nr.O           # draw_opening_data_hora
    $ws.Cells.Item($r, 16).Value = $nr.P           # draw_closing_odds
    $ws.Cells.Item($r, 17).Value = $nr.Q           # draw_closing_data_hora
    $ws.Cells.Item($r, 18).Value = $nr.R           # away_opening_odds
    $ws.Cells.Item($r, 19).Value = $nr.S           # away_opening_data_hora
    $ws.Cells.Item($r, 20).Value = $nr.T           # away_closing_odds
    $ws.Cells.Item($r, 21).Value = $nr.U           # away_closing_data_hora
    $ws.Cells.Item($r, 22).Value = $nr.V           # url_partida
}
